# "Generate Report for Handoff"
# The ea31f9a5-... file became ready for a new handoff cycle:
#  - its Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#    on the Overview sheet (both locale columns) and on each locale sheet.
#  - the "Latest Handoff Datetime" for each locale is refreshed to the new
#    handoff timestamp (shared with the other file's row for that locale,
#    since it records the latest handoff run).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(3, 2).Value = "Ready for handoff"
$overview.Cells.Item(3, 3).Value = "Ready for handoff"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Cells.Item(3, 2).Value = "Ready for handoff"
$zhcn.Cells.Item(2, 4).Value = "2016-03-02 15:39:00"
$zhcn.Cells.Item(3, 4).Value = "2016-03-02 15:39:00"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Cells.Item(3, 2).Value = "Ready for handoff"
$dede.Cells.Item(2, 4).Value = "2016-03-02 15:39:11"
$dede.Cells.Item(3, 4).Value = "2016-03-02 15:39:11"
